$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BECF-pre-ret")

$values = @{
    2  = 0.65404
    3  = 0.67286
    4  = 0.84329
    5  = 0.04674
    6  = 0.33692
    7  = 0.25288
    8  = 0.17033
    9  = 0.66191
    10 = 0.64695
    11 = 0.56819
    12 = 0.33765
    13 = 0.66191
    15 = 0.56819
    16 = 0.56819
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Range("B$row`:AJ$row").Value = $val
}
